$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 628.8570999999999
$ws.Range("I33").Value = 683.6667
$ws.Range("K33").Value = 683.6667
$ws.Range("M33").Value = -454.6667
$ws.Range("H98").Value = 5682390
$ws.Range("I98").Value = 6579504
$ws.Range("K98").Value = 6579504
$ws.Range("M98").Value = -6578006
$ws.Range("H103").Value = 31251062
$ws.Range("I103").Value = 754.1429000000001
$ws.Range("J103").Value = 55556856
$ws.Range("K103").Value = 2262.4287
$ws.Range("L103").Value = 166670568
$ws.Range("M103").Value = -1676.4287
$ws.Range("N103").Value = -166671740
$ws.Range("H122").Value = 5682390
$ws.Range("I122").Value = 6579504
$ws.Range("K122").Value = 19738512
$ws.Range("M122").Value = -19736062
$ws.Range("H137").Value = 14287987
$ws.Range("I137").Value = 19232742
$ws.Range("J137").Value = 3139.2222
$ws.Range("K137").Value = 57698226
$ws.Range("L137").Value = 9417.6666
$ws.Range("M137").Value = -57695676
$ws.Range("N137").Value = -14517.6666
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2660.8
$ws.Range("I2").Value = 2265.1538
$ws.Range("K2").Value = 2265.1538
$ws.Range("M2").Value = -2152.1538
$ws.Range("H4").Value = 4103.3
$ws.Range("J4").Value = 2999
$ws.Range("L4").Value = 2999
$ws.Range("N4").Value = -3231
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -188
$ws.Range("H11").Value = 1166666.6
$ws.Range("I11").Value = 1166666.6
$ws.Range("K11").Value = 1166666.6
$ws.Range("M11").Value = -1166522.6
$ws.Range("H12").Value = 4000
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 5200
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 5200
$ws.Range("M12").Value = -827
$ws.Range("N12").Value = -5546
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670
$ws.Range("H27").Value = 105263
$ws.Range("J27").Value = 105263
$ws.Range("L27").Value = 105263
$ws.Range("N27").Value = -105631
$ws.Range("H39").Value = 5507.5
$ws.Range("I39").Value = 5507.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 5507.5
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4987.5
$ws.Range("N39").ClearContents()
$ws.Range("H58").Value = 89499.5
$ws.Range("I58").Value = 89000
$ws.Range("J58").Value = 89999
$ws.Range("K58").Value = 89000
$ws.Range("L58").Value = 89999
$ws.Range("M58").Value = -88570
$ws.Range("N58").Value = -90859
$ws.Range("H102").Value = 3822.375
$ws.Range("I102").Value = 2906.4443
$ws.Range("K102").Value = 2906.4443
$ws.Range("M102").Value = -1284.4443
$ws.Range("H116").Value = 2660.8
$ws.Range("I116").Value = 2265.1538
$ws.Range("K116").Value = 2265.1538
$ws.Range("M116").Value = 28.84619999999995
$ws.Range("H122").Value = 4395.9
$ws.Range("I122").Value = 3884.3333
$ws.Range("K122").Value = 11652.9999
$ws.Range("M122").Value = -9202.999899999999
$ws.Range("H132").Value = 3487.8
$ws.Range("I132").Value = 3593.5833
$ws.Range("K132").Value = 10780.7499
$ws.Range("M132").Value = -8250.749899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2660.8
$ws.Range("I3").Value = 2265.1538
$ws.Range("K3").Value = 2265.1538
$ws.Range("M3").Value = -2151.1538
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -185
$ws.Range("H22").Value = 5049.75
$ws.Range("I22").Value = 3599.3333
$ws.Range("K22").Value = 3599.3333
$ws.Range("M22").Value = -3426.3333
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7697.2856
$ws.Range("I70").Value = 6935.625
$ws.Range("J70").Value = 8338.684999999999
$ws.Range("K70").Value = 6935.625
$ws.Range("L70").Value = 8338.684999999999
$ws.Range("M70").Value = -6665.625
$ws.Range("N70").Value = -8878.684999999999
$ws.Range("H73").Value = 7697.2856
$ws.Range("I73").Value = 6935.625
$ws.Range("J73").Value = 8338.684999999999
$ws.Range("K73").Value = 6935.625
$ws.Range("L73").Value = 8338.684999999999
$ws.Range("M73").Value = -5999.625
$ws.Range("N73").Value = -10210.685
$ws.Range("H80").Value = 54548210
$ws.Range("I80").Value = 2663.625
$ws.Range("J80").Value = 200002990
$ws.Range("K80").Value = 2663.625
$ws.Range("L80").Value = 200002990
$ws.Range("M80").Value = -1665.625
$ws.Range("N80").Value = -200004986
$ws.Range("H83").Value = 54548210
$ws.Range("I83").Value = 2663.625
$ws.Range("J83").Value = 200002990
$ws.Range("K83").Value = 13318.125
$ws.Range("L83").Value = 1000014950
$ws.Range("M83").Value = -8326.125
$ws.Range("N83").Value = -1000024934
$ws.Range("H102").Value = 1976.5238
$ws.Range("I102").Value = 1921.421
$ws.Range("K102").Value = 1921.421
$ws.Range("M102").Value = -299.421
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1829.6
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 1812
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 1812
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -2188
$ws.Range("H64").Value = 59999.668
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59999.668
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59999.668
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -60449.668
$ws.Range("H67").Value = 59999.668
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59999.668
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59999.668
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -61559.668
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24450
$ws.Range("H63").Value = 29499.5
$ws.Range("J63").Value = 29499.5
$ws.Range("L63").Value = 29499.5
$ws.Range("N63").Value = -30747.5
$ws.Range("H66").Value = 29499.5
$ws.Range("J66").Value = 29499.5
$ws.Range("L66").Value = 88498.5
$ws.Range("N66").Value = -94738.5
$ws.Range("H107").Value = 3940.5806
$ws.Range("J107").Value = 5474.6
$ws.Range("L107").Value = 16423.8
$ws.Range("N107").Value = -20263.8
$ws.Range("H136").Value = 858.8077
$ws.Range("I136").Value = 461.72974
$ws.Range("J136").Value = 1838.2667
$ws.Range("K136").Value = 1385.18922
$ws.Range("L136").Value = 5514.800099999999
$ws.Range("M136").Value = 1164.81078
$ws.Range("N136").Value = -10614.8001
